$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.872.15'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '1.831.53'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = "'310.69"
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').Value = "'0.4624"
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').Value = "'0.3704"
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('E9').Value = '  -2.32%  '
$ws.Range('D10').Value = "'0.8775"
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').Value = "'0.07860"
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').Value = "'19.61"
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').Value = '1.834.65'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = "'5.328"
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('D16').Value = "'86.96"
$ws.Range('E16').Value = '  -5.63%  '
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = "'0.000008735"
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '26.910.32'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('D22').Value = "'4.991"
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').Value = "'1.977"
$ws.Range('D25').Value = "'150.89"
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').Value = "'18.24"
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('D27').Value = "'1.971"
$ws.Range('E27').Value = '  -4.91%  '
$ws.Range('D28').Value = "'113.57"
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = "'4.935"
$ws.Range('E29').Value = '  -3.75%  '
$ws.Range('D30').Value = "'0.08824"
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').Value = "'3.131"
$ws.Range('E31').Value = '  +3.48%  '
$ws.Range('D32').Value = "'0.7569"
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = "'4.463"
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  -2.48%  '
$ws.Range('D35').Value = "'2.600"
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('D37').Value = "'0.01934"
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('D38').Value = "'2.926"
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('D40').Value = "'6.903"
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('D41').Value = "'0.4976"
$ws.Range('E41').Value = '  -3.62%  '
$ws.Range('D42').Value = "'0.1597"
$ws.Range('E42').Value = '  -2.89%  '
$ws.Range('D43').Value = "'8.349"
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = "'0.4676"
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('D45').Value = "'1.007"
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').Value = "'10.12"
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('D47').Value = "'102.39"
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('E48').Value = '  -2.43%  '
$ws.Range('D49').Value = "'0.06100"
$ws.Range('E49').Value = '  -2.20%  '
$ws.Range('D50').Value = "'64.49"
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').Value = "'36.49"
$ws.Range('E51').Value = '  -1.58%  '
